$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'25.910.61"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  -0.16%  '
$ws.Cells.Item(3, 4).Value = "'1.639.30"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  -0.16%  '
$ws.Cells.Item(4, 4).Value = "'1.003"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.10%  '
$ws.Cells.Item(5, 4).Value = "'214.74"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -0.21%  '
$ws.Cells.Item(6, 4).Value = "'0.5095"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +0.99%  '
$ws.Cells.Item(7, 4).Value = "'1.004"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -0.01%  '
$ws.Cells.Item(8, 4).Value = "'0.2561"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -0.42%  '
$ws.Cells.Item(9, 4).Value = "'0.06369"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -0.70%  '
$ws.Cells.Item(10, 4).Value = "'19.52"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -0.44%  '
$ws.Cells.Item(11, 4).Value = "'0.07746"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -0.44%  '
$ws.Cells.Item(12, 4).Value = "'1.655.54"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +2.03%  '
$ws.Cells.Item(13, 4).Value = "'4.279"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +0.31%  '
$ws.Cells.Item(14, 4).Value = "'0.5437"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -0.08%  '
$ws.Cells.Item(15, 4).Value = "'0.0₅7801"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -1.64%  '
$ws.Cells.Item(16, 4).Value = "'64.18"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -0.40%  '
$ws.Cells.Item(17, 4).Value = "'25.949.74"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -0.15%  '
$ws.Cells.Item(18, 5).Value = '  -0.15%  '
$ws.Cells.Item(19, 4).Value = "'196.63"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -1.72%  '
$ws.Cells.Item(20, 4).Value = "'4.421"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +0.94%  '
$ws.Cells.Item(21, 4).Value = "'9.927"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +0.25%  '
$ws.Cells.Item(22, 4).Value = "'6.035"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +0.98%  '
$ws.Cells.Item(23, 4).Value = "'1.006"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +0.15%  '
$ws.Cells.Item(24, 4).Value = "'1.879"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -0.63%  '
$ws.Cells.Item(25, 4).Value = "'141.60"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +0.65%  '
$ws.Cells.Item(26, 4).Value = "'0.1204"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +6.10%  '
$ws.Cells.Item(27, 4).Value = "'6.847"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +0.33%  '
$ws.Cells.Item(28, 4).Value = "'15.67"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  -0.24%  '
$ws.Cells.Item(29, 5).Value = '  -0.52%  '
$ws.Cells.Item(30, 4).Value = "'0.04947"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +0.51%  '
$ws.Cells.Item(31, 4).Value = "'3.248"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -0.45%  '
$ws.Cells.Item(32, 4).Value = "'3.171"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -1.24%  '
$ws.Cells.Item(33, 4).Value = "'1.532"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(34, 4).Value = "'2.370"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -0.13%  '
$ws.Cells.Item(35, 4).Value = "'0.8924"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +0.00%  '
$ws.Cells.Item(36, 4).Value = "'2.578"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -1.05%  '
$ws.Cells.Item(37, 4).Value = "'1.134.40"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -1.82%  '
$ws.Cells.Item(38, 4).Value = "'0.5419"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -2.66%  '
$ws.Cells.Item(39, 4).Value = "'0.01552"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -1.13%  '
$ws.Cells.Item(40, 5).Value = '  -0.05%  '
$ws.Cells.Item(41, 5).Value = '  -0.57%  '
$ws.Cells.Item(42, 4).Value = "'5.561"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -2.75%  '
$ws.Cells.Item(43, 4).Value = "'0.8143"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +0.43%  '
$ws.Cells.Item(44, 4).Value = "'99.46"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -0.29%  '
$ws.Cells.Item(45, 4).Value = "'0.0₈123"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +2.44%  '
$ws.Cells.Item(46, 4).Value = "'1.776.23"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -0.21%  '
$ws.Cells.Item(47, 4).Value = "'0.4537"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +0.33%  '
$ws.Cells.Item(48, 4).Value = "'1.004"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -0.06%  '
$ws.Cells.Item(49, 4).Value = "'54.71"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -0.16%  '
$ws.Cells.Item(50, 4).Value = "'0.05063"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +0.24%  '
$ws.Cells.Item(51, 5).Value = '  +0.14%  '
